# Update "Countries & provincias Spain" data
# - refresh the "last updated" timestamp
# - update a handful of per-country stat rows (some of which also
#   change relative order vs. their neighbour because the updated
#   "Casos totales" value changes their rank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row --------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 06:05"

# --- Row 8: Brasil (unchanged rank, value refresh) --------------------------
$ws.Range("A8").Value = "Brasil"
$ws.Range("B8").Value = 241080
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 94122
$ws.Range("E8").Value = 130836
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 16122

# --- Row 23: Pakistan (unchanged rank, value refresh) -----------------------
$ws.Range("A23").Value = "Pakistan"
$ws.Range("B23").Value = 42125
$ws.Range("C23").Value = 1974
$ws.Range("D23").Value = 11922
$ws.Range("E23").Value = 29300
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 903

# --- Rows 60-61: Kazajistan overtakes Finlandia -----------------------------
$ws.Range("A60").Value = "Kazajistan"
$ws.Range("B60").Value = 6440
$ws.Range("C60").Value = 283
$ws.Range("D60").Value = 3256
$ws.Range("E60").Value = 3150
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 34

$ws.Range("A61").Value = "Finlandia"
$ws.Range("B61").Value = 6347
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 5000
$ws.Range("E61").Value = 1049
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 298

# --- Rows 143-145: Nepal overtakes Madagascar and Togo ----------------------
$ws.Range("A143").Value = "Nepal"
$ws.Range("B143").Value = 304
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 36
$ws.Range("E143").Value = 266
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 2

$ws.Range("A144").Value = "Madagascar"
$ws.Range("B144").Value = 304
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 114
$ws.Range("E144").Value = 189
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1

$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 301
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 104
$ws.Range("E145").Value = 186
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 11
